$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

$ws.Range("A$row").Value = 112209549
$ws.Range("B$row").Value = 60151
$ws.Range("C$row").Value = "Ovaliderad"
$ws.Range("D$row").Value = "NT"
$ws.Range("E$row").Value = 100452
$ws.Range("F$row").Value = "Tvåtandad spolsnäcka"
$ws.Range("G$row").Value = "Alinda biplicata"
$ws.Range("H$row").Value = "(Montagu, 1803)"
$ws.Range("P$row").Value = "Flugmötesskogen , Srm"
$ws.Range("Q$row").Value = 580550
$ws.Range("R$row").Value = 6579320
$ws.Range("S$row").Value = 5
$ws.Range("T$row").Value = "Södermanland"
$ws.Range("U$row").Value = "Eskilstuna"
$ws.Range("V$row").Value = "Södermanland"
$ws.Range("W$row").Value = "Eskilstuna"

$ws.Range("Y$row").NumberFormat = "@"
$ws.Range("Y$row").Value = "2023-09-19"
$ws.Range("Z$row").NumberFormat = "@"
$ws.Range("Z$row").Value = "11:02"
$ws.Range("AA$row").NumberFormat = "@"
$ws.Range("AA$row").Value = "2023-09-19"
$ws.Range("AB$row").NumberFormat = "@"
$ws.Range("AB$row").Value = "11:02"

$ws.Range("AC$row").Value = "Åt på en hasselticka"
$ws.Range("AD$row").Value = $false
$ws.Range("AE$row").Value = $false
$ws.Range("AG$row").Value = $false
$ws.Range("AW$row").Value = "Ella Axelsson Elfving"
$ws.Range("AX$row").Value = "Ella Axelsson Elfving"
